$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh1 = $s.Shapes.AddShape(1, 4191814/12700, 2574518/12700, 432000/12700, 432000/12700)
Write-Host "before: W=$($sh1.Width)"
$factor = (118727/12700) / $sh1.Width
Write-Host "factor=$factor"
$sh1.ScaleWidth($factor, 0)
Write-Host "after: W=$($sh1.Width)"
